$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert 3 new rows before row 125. This shifts the previous rows 125-189
# down to 128-192, preserving all their data and formatting (including the
# date-formatted style used in column D).
$ws.Rows("125:127").Insert()

# Populate the 3 newly inserted rows (125-127) with the new data records.
# Row 125: Especial
$ws.Cells.Item(125, 1).Value2 = 11
$ws.Cells.Item(125, 2).Value2 = "Vega Monumental Concepción"
$ws.Cells.Item(125, 3).Value2 = "Bíobío"
$ws.Cells.Item(125, 4).Value2 = 44777
$ws.Cells.Item(125, 5).Value2 = 8
$ws.Cells.Item(125, 6).Value2 = "Fruta"
$ws.Cells.Item(125, 7).Value2 = 100101
$ws.Cells.Item(125, 8).Value2 = "Berries"
$ws.Cells.Item(125, 9).Value2 = 100101007
$ws.Cells.Item(125, 10).Value2 = "Kiwi"
$ws.Cells.Item(125, 11).Value2 = "Hayward"
$ws.Cells.Item(125, 12).Value2 = "Especial"
$ws.Cells.Item(125, 13).Value2 = 50
$ws.Cells.Item(125, 14).Value2 = 9000
$ws.Cells.Item(125, 15).Value2 = 9000
$ws.Cells.Item(125, 16).Value2 = 9000
$ws.Cells.Item(125, 17).Value2 = "$/bandeja 18 kilos"
$ws.Cells.Item(125, 18).Value2 = "Región de O'Higgins"
$ws.Cells.Item(125, 19).Value2 = 500
$ws.Cells.Item(125, 20).Value2 = 18

# Row 126: Primera
$ws.Cells.Item(126, 1).Value2 = 11
$ws.Cells.Item(126, 2).Value2 = "Vega Monumental Concepción"
$ws.Cells.Item(126, 3).Value2 = "Bíobío"
$ws.Cells.Item(126, 4).Value2 = 44777
$ws.Cells.Item(126, 5).Value2 = 8
$ws.Cells.Item(126, 6).Value2 = "Fruta"
$ws.Cells.Item(126, 7).Value2 = 100101
$ws.Cells.Item(126, 8).Value2 = "Berries"
$ws.Cells.Item(126, 9).Value2 = 100101007
$ws.Cells.Item(126, 10).Value2 = "Kiwi"
$ws.Cells.Item(126, 11).Value2 = "Hayward"
$ws.Cells.Item(126, 12).Value2 = "Primera"
$ws.Cells.Item(126, 13).Value2 = 50
$ws.Cells.Item(126, 14).Value2 = 8000
$ws.Cells.Item(126, 15).Value2 = 8000
$ws.Cells.Item(126, 16).Value2 = 8000
$ws.Cells.Item(126, 17).Value2 = "$/bandeja 18 kilos"
$ws.Cells.Item(126, 18).Value2 = "Región de O'Higgins"
$ws.Cells.Item(126, 19).Value2 = 444
$ws.Cells.Item(126, 20).Value2 = 18

# Row 127: Segunda
$ws.Cells.Item(127, 1).Value2 = 11
$ws.Cells.Item(127, 2).Value2 = "Vega Monumental Concepción"
$ws.Cells.Item(127, 3).Value2 = "Bíobío"
$ws.Cells.Item(127, 4).Value2 = 44777
$ws.Cells.Item(127, 5).Value2 = 8
$ws.Cells.Item(127, 6).Value2 = "Fruta"
$ws.Cells.Item(127, 7).Value2 = 100101
$ws.Cells.Item(127, 8).Value2 = "Berries"
$ws.Cells.Item(127, 9).Value2 = 100101007
$ws.Cells.Item(127, 10).Value2 = "Kiwi"
$ws.Cells.Item(127, 11).Value2 = "Hayward"
$ws.Cells.Item(127, 12).Value2 = "Segunda"
$ws.Cells.Item(127, 13).Value2 = 50
$ws.Cells.Item(127, 14).Value2 = 7000
$ws.Cells.Item(127, 15).Value2 = 7000
$ws.Cells.Item(127, 16).Value2 = 7000
$ws.Cells.Item(127, 17).Value2 = "$/bandeja 18 kilos"
$ws.Cells.Item(127, 18).Value2 = "Región de O'Higgins"
$ws.Cells.Item(127, 19).Value2 = 389
$ws.Cells.Item(127, 20).Value2 = 18
